$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 3.7
$ws.Range("L2").Value = 4.75
$ws.Range("N2").Value = 4.75
$ws.Range("S2").Value = 1.78
$ws.Range("T2").Value = 2.03
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 29
$ws.Range("AC2").Value = 4.75
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 126
$ws.Range("AG2").Value = 6.5
$ws.Range("AO2").Value = 17
$ws.Range("AP2").Value = 41
$ws.Range("AR2").Value = 126
$ws.Range("AX2").Value = 26
$ws.Range("AZ2").Value = 101

# Row 5
$ws.Range("G5").Value = 1.75
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 2.4
$ws.Range("L5").Value = 5
$ws.Range("Y5").Value = 9
$ws.Range("AC5").Value = 8.5
$ws.Range("AM5").Value = 700
$ws.Range("AO5").Value = 9.5
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 34
$ws.Range("AX5").Value = 26

# Row 6
$ws.Range("G6").Value = 1.36
$ws.Range("H6").Value = 4.8
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 1.75
$ws.Range("K6").Value = 2.57
$ws.Range("L6").Value = 6.5
$ws.Range("N6").Value = 9.5
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4.8
$ws.Range("Q6").Value = 1.45
$ws.Range("R6").Value = 2.55
$ws.Range("S6").Value = 1.26
$ws.Range("T6").Value = 3.55
$ws.Range("V6").Value = 2.05
$ws.Range("W6").Value = 9.5
$ws.Range("X6").Value = 8
$ws.Range("Z6").Value = 9.75
$ws.Range("AA6").Value = 10
$ws.Range("AB6").Value = 20
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 9.75
$ws.Range("AE6").Value = 17
$ws.Range("AG6").Value = 24
$ws.Range("AH6").Value = 55
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 175
$ws.Range("AK6").Value = 75
$ws.Range("AL6").Value = 60
$ws.Range("AN6").Value = 3.5
$ws.Range("AO6").Value = 6
$ws.Range("AP6").Value = 13
$ws.Range("AQ6").Value = 15
$ws.Range("AR6").Value = 32
$ws.Range("AS6").Value = 120
$ws.Range("AT6").Value = 3.55
$ws.Range("AU6").Value = 7.6
$ws.Range("AW6").Value = 9
$ws.Range("AX6").Value = 40
$ws.Range("AY6").Value = 35
$ws.Range("AZ6").Value = 250
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 400

# Row 7
$ws.Range("G7").Value = 2.27
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 2.92
$ws.Range("J7").Value = 2.87
$ws.Range("L7").Value = 3.55
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 3.8
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 2.07
$ws.Range("T7").Value = 2.72
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.32
$ws.Range("W7").Value = 10.25
$ws.Range("X7").Value = 13.5
$ws.Range("Y7").Value = 8.75
$ws.Range("AA7").Value = 16.5
$ws.Range("AB7").Value = 21
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = 6.6
$ws.Range("AE7").Value = 11.25
$ws.Range("AF7").Value = 40
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 17
$ws.Range("AK7").Value = 23
$ws.Range("AL7").Value = 26
$ws.Range("AM7").Value = 250
$ws.Range("AO7").Value = 12
$ws.Range("AP7").Value = 18.5
$ws.Range("AR7").Value = 75
$ws.Range("AS7").Value = 200
$ws.Range("AT7").Value = 2.72
$ws.Range("AW7").Value = 5

# Row 8
$ws.Range("G8").Value = 1.88
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 3.65
$ws.Range("J8").Value = 2.47
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 8.25
$ws.Range("O8").Value = 1.23
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.07
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 2.92
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.15
$ws.Range("W8").Value = 8.25
$ws.Range("X8").Value = 9.75
$ws.Range("AA8").Value = 14
$ws.Range("AB8").Value = 22
$ws.Range("AC8").Value = 8.25
$ws.Range("AD8").Value = 7.1
$ws.Range("AF8").Value = 50
$ws.Range("AG8").Value = 13
$ws.Range("AH8").Value = 22
$ws.Range("AM8").Value = 350
$ws.Range("AN8").Value = 3.9
$ws.Range("AO8").Value = 9.5
$ws.Range("AP8").Value = 17
$ws.Range("AR8").Value = 60
$ws.Range("AT8").Value = 2.92
$ws.Range("AU8").Value = 6.9
$ws.Range("AV8").Value = 55
$ws.Range("AX8").Value = 19.5
$ws.Range("AY8").Value = 24
$ws.Range("BB8").Value = 250

# Row 9
$ws.Range("G9").Value = 2.72
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 2.47
$ws.Range("N9").Value = 7.3
$ws.Range("U9").Value = 1.65
$ws.Range("V9").Value = 2.1
$ws.Range("W9").Value = 9.25
$ws.Range("X9").Value = 14.5
$ws.Range("Y9").Value = 9.75
$ws.Range("AC9").Value = 7.3
$ws.Range("AG9").Value = 8.75
$ws.Range("AH9").Value = 13
$ws.Range("AJ9").Value = 27
$ws.Range("AK9").Value = 19.5
$ws.Range("AL9").Value = 27
$ws.Range("AN9").Value = 4.7
$ws.Range("AQ9").Value = 70
$ws.Range("AS9").Value = 300
$ws.Range("AU9").Value = 6.9
$ws.Range("AW9").Value = 4.45
$ws.Range("AZ9").Value = 55

# Row 10
$ws.Range("G10").Value = 1.38
$ws.Range("H10").Value = 4.4
$ws.Range("I10").Value = 7.7
$ws.Range("K10").Value = 2.27
$ws.Range("L10").Value = 7
$ws.Range("N10").Value = 7.7
$ws.Range("S10").Value = 1.39
$ws.Range("T10").Value = 2.77
$ws.Range("W10").Value = 6.2
$ws.Range("X10").Value = 6
$ws.Range("AC10").Value = 7.7
$ws.Range("AG10").Value = 18
$ws.Range("AH10").Value = 50
$ws.Range("AJ10").Value = 200
$ws.Range("AK10").Value = 100
$ws.Range("AQ10").Value = 19
$ws.Range("AT10").Value = 2.77
$ws.Range("AW10").Value = 8.5

# Row 11
$ws.Range("G11").Value = 1.24
$ws.Range("I11").Value = 11
$ws.Range("K11").Value = 2.72
$ws.Range("L11").Value = 7.9
$ws.Range("P11").Value = 5.4
$ws.Range("Q11").Value = 1.38
$ws.Range("R11").Value = 2.82
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.95
$ws.Range("AA11").Value = 9.75
$ws.Range("AE11").Value = 19.5
$ws.Range("AF11").Value = 65
$ws.Range("AG11").Value = 40
$ws.Range("AH11").Value = 100
$ws.Range("AI11").Value = 32
$ws.Range("AJ11").Value = 350
$ws.Range("AK11").Value = 120
$ws.Range("AS11").Value = 120
$ws.Range("AW11").Value = 11.25
$ws.Range("AX11").Value = 55
$ws.Range("AZ11").Value = 400
$ws.Range("BA11").Value = 300

